$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.035.31"
$ws.Range("E2").Value = "  -4.90%  "

$ws.Range("D3").Value = "3.330.18"
$ws.Range("E3").Value = "  -5.94%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.58"
$ws.Range("E5").Value = "  -4.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.62"
$ws.Range("E6").Value = "  -7.14%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -3.50%  "

$ws.Range("D9").Value = "3.323.84"
$ws.Range("E9").Value = "  -5.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.184"
$ws.Range("E10").Value = "  -10.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.583"
$ws.Range("E11").Value = "  -7.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.15"
$ws.Range("E12").Value = "  -8.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("E13").Value = "  -7.70%  "

$ws.Range("D14").Value = "3.866.58"
$ws.Range("E14").Value = "  -5.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.52"
$ws.Range("E15").Value = "  -7.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "600.59"
$ws.Range("E16").Value = "  -9.78%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.202.77"
$ws.Range("E17").Value = "  -4.77%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.92"
$ws.Range("E18").Value = "  -2.88%  "

$ws.Range("D19").Value = "3.335.72"
$ws.Range("E19").Value = "  -5.35%  "

$ws.Range("E20").Value = "  -3.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.34"
$ws.Range("E21").Value = "  -9.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.899"
$ws.Range("E22").Value = "  -6.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.81"
$ws.Range("E23").Value = "  -7.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.08"
$ws.Range("E24").Value = "  -4.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.04"
$ws.Range("E25").Value = "  -4.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.01"
$ws.Range("E26").Value = "  -8.19%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.64"
$ws.Range("E28").Value = "  -8.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.23"
$ws.Range("E29").Value = "  -9.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.67"
$ws.Range("E30").Value = "  -9.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.58"
$ws.Range("E31").Value = "  -7.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.22"
$ws.Range("E32").Value = "  -7.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.70"
$ws.Range("E33").Value = "  -15.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.97"
$ws.Range("E34").Value = "  -6.60%  "

$ws.Range("E35").Value = "  -6.26%  "

$ws.Range("D36").Value = "3.775.93"
$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.84"
$ws.Range("E37").Value = "  -6.38%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "531.31"
$ws.Range("E38").Value = "  +6.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.45"
$ws.Range("E40").Value = "  -7.00%  "

$ws.Range("D41").Value = "0.0₃0710"
$ws.Range("E41").Value = "  -12.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  -8.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.124"
$ws.Range("E43").Value = "  -7.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.338"
$ws.Range("E44").Value = "  -8.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.67"
$ws.Range("E45").Value = "  -8.69%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.25"
$ws.Range("E46").Value = "  -3.23%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0412"
$ws.Range("E47").Value = "  -8.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.00"
$ws.Range("E48").Value = "  +10.65%  "

$ws.Range("E49").Value = "  -5.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.60"
$ws.Range("E50").Value = "  -8.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.10%  "
